$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-12-26 Thursday" "2024-12-27 Friday"

Replace-Text "31×72=" "39×54="
Replace-Text "96×17=" "43×15="
Replace-Text "32×96=" "36×85="
Replace-Text "64×36=" "68×53="
Replace-Text "69×80=" "42×89="
Replace-Text "74×43=" "80×61="
Replace-Text "51×22=" "96×65="
Replace-Text "17×14=" "53×72="
Replace-Text "65×99=" "95×75="
Replace-Text "44×51=" "52×52="
Replace-Text "35×86=" "91×82="
Replace-Text "14×36=" "74×25="
Replace-Text "69×81=" "21×87="
Replace-Text "34×98=" "23×59="
Replace-Text "98×60=" "86×71="
Replace-Text "85×95=" "95×28="
Replace-Text "21×92=" "97×40="
Replace-Text "45×23=" "15×58="
Replace-Text "19×17=" "78×83="
Replace-Text "50×41=" "28×12="
Replace-Text "43×99=" "75×17="
Replace-Text "95×44=" "74×70="
Replace-Text "63×42=" "87×65="
Replace-Text "53×30=" "51×64="
Replace-Text "37×42=" "34×21="
